# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# OFF sheet - row 3 (R)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 391
$wsOff.Range("C3").Value = 284
$wsOff.Range("D3").Value = 105
$wsOff.Range("E3").Value = 48

# DEF sheet - row 3 (R)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 444
$wsDef.Range("C3").Value = 309
$wsDef.Range("D3").Value = 93
$wsDef.Range("E3").Value = 44
